$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, formatted like the other headers (copy style from E1)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "Modelo"

$modelo = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       LinearRegression())]),`n                                            param_grid={'model__fit_intercept': [True,`n                                                                                 False]},`n                                            scoring='neg_mean_squared_error'))"

$ws.Range("F2").Value = $modelo
$ws.Range("F3").Value = $modelo
$ws.Range("F4").Value = $modelo
$ws.Range("F5").Value = $modelo
